$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $origStyle = $Range.Style
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '66.801.70'
Set-TextValue $ws.Range('E2') '  -3.99%  '
Set-TextValue $ws.Range('D3') '3.459.91'
Set-TextValue $ws.Range('E3') '  -4.27%  '
Set-TextValue $ws.Range('E4') '  +0.09%  '
Set-TextValue $ws.Range('D5') '602.81'
Set-TextValue $ws.Range('E5') '  -4.38%  '
Set-TextValue $ws.Range('D6') '147.32'
Set-TextValue $ws.Range('E6') '  -7.14%  '
Set-TextValue $ws.Range('D7') '3.460.69'
Set-TextValue $ws.Range('E7') '  -4.22%  '
Set-TextValue $ws.Range('E8') '  +0.03%  '
Set-TextValue $ws.Range('E9') '  -2.30%  '
Set-TextValue $ws.Range('E10') '  -5.04%  '
Set-TextValue $ws.Range('D11') '7.47'
Set-TextValue $ws.Range('E11') '  -0.49%  '
Set-TextValue $ws.Range('D12') '0.424'
Set-TextValue $ws.Range('E12') '  -4.25%  '
Set-TextValue $ws.Range('E13') '  -5.76%  '
Set-TextValue $ws.Range('D14') '31.63'
Set-TextValue $ws.Range('E14') '  -6.23%  '
Set-TextValue $ws.Range('D15') '4.051.84'
Set-TextValue $ws.Range('E15') '  -4.07%  '
Set-TextValue $ws.Range('D16') '3.460.86'
Set-TextValue $ws.Range('E16') '  -4.18%  '
Set-TextValue $ws.Range('D17') '66.926.50'
Set-TextValue $ws.Range('E17') '  -3.68%  '
Set-TextValue $ws.Range('E18') '  -0.91%  '
Set-TextValue $ws.Range('D19') '6.44'
Set-TextValue $ws.Range('E19') '  -4.62%  '
Set-TextValue $ws.Range('E20') '  -5.55%  '
Set-TextValue $ws.Range('D21') '10.00'
Set-TextValue $ws.Range('E21') '  -2.89%  '
Set-TextValue $ws.Range('D22') '439.89'
Set-TextValue $ws.Range('E22') '  -4.84%  '
Set-TextValue $ws.Range('D23') '0.609'
Set-TextValue $ws.Range('E23') '  -5.98%  '
Set-TextValue $ws.Range('D24') '78.75'
Set-TextValue $ws.Range('E24') '  -0.32%  '
Set-TextValue $ws.Range('E25') '  -0.03%  '
Set-TextValue $ws.Range('D26') '3.605.36'
Set-TextValue $ws.Range('E26') '  -4.04%  '
Set-TextValue $ws.Range('E27') '  -10.59%  '
Set-TextValue $ws.Range('D28') '9.90'
Set-TextValue $ws.Range('E28') '  -8.07%  '
Set-TextValue $ws.Range('D29') '8.40'
Set-TextValue $ws.Range('E29') '  -10.37%  '
Set-TextValue $ws.Range('E30') '  -6.50%  '
Set-TextValue $ws.Range('E31') '  -7.23%  '
Set-TextValue $ws.Range('E32') '  -3.13%  '
Set-TextValue $ws.Range('D33') '1.00'
Set-TextValue $ws.Range('E33') '  +0.06%  '
Set-TextValue $ws.Range('D34') '25.37'
Set-TextValue $ws.Range('E34') '  -4.58%  '
Set-TextValue $ws.Range('D35') '6.07'
Set-TextValue $ws.Range('E35') '  -7.60%  '
Set-TextValue $ws.Range('D36') '3.458.53'
Set-TextValue $ws.Range('E36') '  -4.19%  '
Set-TextValue $ws.Range('E37') '  -7.74%  '
Set-TextValue $ws.Range('B38') 'USDe'
Set-TextValue $ws.Range('C38') 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range('D38') '1.00'
Set-TextValue $ws.Range('E38') '  +0.01%  '
Set-TextValue $ws.Range('B39') 'Aptos'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D39') '7.91'
Set-TextValue $ws.Range('E39') '  -6.84%  '
Set-TextValue $ws.Range('D40') '0.999'
Set-TextValue $ws.Range('E40') '  -0.02%  '
Set-TextValue $ws.Range('D41') '2.17'
Set-TextValue $ws.Range('E41') '  -10.37%  '
Set-TextValue $ws.Range('B42') 'Hedera'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D42') '0.0888'
Set-TextValue $ws.Range('E42') '  -4.36%  '
Set-TextValue $ws.Range('B43') 'Monero'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D43') '171.86'
Set-TextValue $ws.Range('E43') '  -3.97%  '
Set-TextValue $ws.Range('E44') '  -5.52%  '
Set-TextValue $ws.Range('D45') '0.883'
Set-TextValue $ws.Range('E45') '  -3.44%  '
Set-TextValue $ws.Range('D46') '28.99'
Set-TextValue $ws.Range('E46') '  -9.14%  '
Set-TextValue $ws.Range('D47') '45.93'
Set-TextValue $ws.Range('E47') '  -0.05%  '
Set-TextValue $ws.Range('E48') '  -10.92%  '
Set-TextValue $ws.Range('E49') '  -4.81%  '
Set-TextValue $ws.Range('D50') '2.45'
Set-TextValue $ws.Range('E50') '  -11.29%  '
Set-TextValue $ws.Range('D51') '0.983'
Set-TextValue $ws.Range('E51') '  -5.76%  '
